# Add season record (Wins / Losses / Ties) columns to the roster sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): AD1=Wins, AE1=Losses, AF1=Ties ---
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the header formatting used by the rest of row 1 (bold, centered, bordered).
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Re-apply the header text (PasteSpecial(xlPasteFormats) only copies formatting,
# but make sure values are still correct).
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# --- Data rows (2-50): every team had the same 2011 record: 82 wins, 79 losses, 0 ties ---
for ($r = 2; $r -le 50; $r++) {
    $ws.Cells.Item($r, 30).Value = 82
    $ws.Cells.Item($r, 31).Value = 79
    $ws.Cells.Item($r, 32).Value = 0
}

Write-Host "Added Wins/Losses/Ties columns (AD:AF) for rows 1-50"
